# Rename the data-dictionary style column headers on the "Covid-19" sheet
# from descriptive labels to snake_case variable names (matches the
# "add_shelter_length" upload).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19")

$ws.Range("B1").Value = "total_cases"
$ws.Range("C1").Value = "total_death"
$ws.Range("D1").Value = "death_100k"
$ws.Range("E1").Value = "cases_last_7days "
$ws.Range("F1").Value = "rate_per_100k"
$ws.Range("G1").Value = "total_test_results"
$ws.Range("H1").Value = "state_of_emergency"
$ws.Range("I1").Value = "stay_ at_ home_ shelter_in_ place"
$ws.Range("J1").Value = "end_relax_stay_at_home_shelter_in_place"
$ws.Range("K1").Value = "length_ shelter_in_place"
$ws.Range("L1").Value = "closed_nonessential_businesses"
$ws.Range("M1").Value = "began_reopen_businesses_statewide"
$ws.Range("N1").Value = "mandate_facemask_use_ employees _ public"
$ws.Range("O1").Value = "weekly_unemployment_insurance_max"
$ws.Range("P1").Value = "population_density_per_sq_mi"
$ws.Range("Q1").Value = "population_2018"
$ws.Range("R1").Value = "percent_living_under_ fed_poverty_line_2018"
$ws.Range("S1").Value = "percent_ risk_ serious_ illness_due_to_covid"
$ws.Range("T1").Value = "all-cause deaths 2018"
$ws.Range("U1").Value = "children_0_18"
$ws.Range("V1").Value = "adults_19_25"
$ws.Range("W1").Value = "adults_26_34"
$ws.Range("X1").Value = "adults_35_54"
$ws.Range("Y1").Value = "adults_55_64"
$ws.Range("AA1").Value = "political_party_governor"

# Cosmetic: scroll/selection state that Excel stored when the file was
# last saved.
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("S1").Select()

$ws1 = $wb.Worksheets.Item("Infromation")
$ws1.Application.ActiveWindow.ScrollRow = 8
$ws1.Range("C22:C26").Select()
